$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.712.23"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "3.845.92"
$ws.Range("E3").Value = "  +2.71%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'609.38"
$ws.Range("E5").Value = "  -2.32%  "
$ws.Range("D6").Value = "'174.49"
$ws.Range("E6").Value = "  -3.73%  "
$ws.Range("D7").Value = "3.845.30"
$ws.Range("E7").Value = "  +2.74%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.526"
$ws.Range("E9").Value = "  -1.71%  "
$ws.Range("D10").Value = "'0.166"
$ws.Range("E10").Value = "  -1.67%  "
$ws.Range("D11").Value = "'6.47"
$ws.Range("E11").Value = "  +2.41%  "
$ws.Range("D12").Value = "'0.480"
$ws.Range("E12").Value = "  -1.53%  "
$ws.Range("D13").Value = "'39.84"
$ws.Range("E13").Value = "  -2.57%  "
$ws.Range("E14").Value = "  -2.78%  "
$ws.Range("D15").Value = "4.476.57"
$ws.Range("E15").Value = "  +2.53%  "
$ws.Range("D16").Value = "3.833.88"
$ws.Range("E16").Value = "  +2.35%  "
$ws.Range("D17").Value = "69.759.40"
$ws.Range("E17").Value = "  -0.73%  "
$ws.Range("D18").Value = "'7.45"
$ws.Range("E18").Value = "  -2.31%  "
$ws.Range("E19").Value = "  -3.23%  "
$ws.Range("D20").Value = "'16.59"
$ws.Range("E20").Value = "  -1.51%  "
$ws.Range("D21").Value = "'504.38"
$ws.Range("E21").Value = "  -0.50%  "
$ws.Range("D22").Value = "'9.50"
$ws.Range("E22").Value = "  +1.60%  "
$ws.Range("D23").Value = "'0.738"
$ws.Range("E23").Value = "  +1.55%  "
$ws.Range("B24").Value = "Fetch.AI"
$ws.Range("C24").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D24").Value = "'2.45"
$ws.Range("E24").Value = "  -5.22%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'85.81"
$ws.Range("E25").Value = "  -1.18%  "
$ws.Range("E26").Value = "  +4.37%  "
$ws.Range("E27").Value = "  -4.67%  "
$ws.Range("D28").Value = "'10.41"
$ws.Range("E28").Value = "  -9.49%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").Value = "'2.52"
$ws.Range("E30").Value = "  +0.69%  "
$ws.Range("D31").Value = "'2.98"
$ws.Range("E31").Value = "  +1.20%  "
$ws.Range("D33").Value = "'32.36"
$ws.Range("E33").Value = "  +3.42%  "
$ws.Range("D34").Value = "'0.113"
$ws.Range("E34").Value = "  -2.66%  "
$ws.Range("D35").Value = "'0.998"
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("E36").Value = "  -1.87%  "
$ws.Range("D37").Value = "'6.08"
$ws.Range("E37").Value = "  -1.71%  "
$ws.Range("E38").Value = "  +1.09%  "
$ws.Range("D39").Value = "'485.40"
$ws.Range("E39").Value = "  +11.94%  "
$ws.Range("D40").Value = "'0.335"
$ws.Range("E40").Value = "  -1.82%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'2.05"
$ws.Range("E41").Value = "  -2.50%  "
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").Value = "'49.70"
$ws.Range("E42").Value = "  -1.28%  "
$ws.Range("D43").Value = "'2.98"
$ws.Range("E43").Value = "  +2.23%  "
$ws.Range("D44").Value = "'43.17"
$ws.Range("E44").Value = "  -6.45%  "
$ws.Range("D45").Value = "'8.50"
$ws.Range("E45").Value = "  -2.78%  "
$ws.Range("D46").Value = "2.916.14"
$ws.Range("E46").Value = "  -3.29%  "
$ws.Range("D47").Value = "'0.0360"
$ws.Range("E47").Value = "  -1.51%  "
$ws.Range("D48").Value = "'139.94"
$ws.Range("E48").Value = "  +1.45%  "
$ws.Range("D50").Value = "'26.77"
$ws.Range("E50").Value = "  -2.88%  "
$ws.Range("E51").Value = "  -4.21%  "
